# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (cloned header/data layout from the
# "2021-Q4" sheet) positioned right after "2021-Q4" and before "总计",
# then records its totals as a new first data row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# ---- 1. New "2022-Q1" sheet, inserted right after "2021-Q4" ----
$ws = $wb.Worksheets.Add($null, $q4)
$ws.Name = "2022-Q1"

# Re-fetch "总计" AFTER the insert: inserting a sheet ahead of it shifts its
# position, so a handle grabbed beforehand would resolve to the wrong tab.
$total = $wb.Worksheets.Item("总计")

# Header row (bold, bordered, centered) mirrors the other quarterly sheets.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 1).Font.Bold = $true
$ws.Cells.Item(2, 1).HorizontalAlignment = -4108
$ws.Cells.Item(2, 1).VerticalAlignment = -4160
$ws.Cells.Item(2, 1).Borders.LineStyle = 1

$ws.Range("B2").Value = "'590003"
$ws.Range("C2").Value = "中邮核心优势灵活配置混合"
$ws.Range("D2").Value = "'12.43"
$ws.Range("E2").Value = "'78.81"
$ws.Range("F2").Value = "'4.81"
$ws.Range("G2").Value = "'0.5979"
$ws.Range("H2").Value = 9

# ---- 2. Record the new quarter as the first data row of "总计" ----
$total.Rows(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 1).Font.Bold = $true
$total.Cells.Item(2, 1).HorizontalAlignment = -4108
$total.Cells.Item(2, 1).VerticalAlignment = -4160
$total.Cells.Item(2, 1).Borders.LineStyle = 1

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0.6
$total.Range("B2:D2").Font.Bold = $false
$total.Range("B2:D2").Borders.LineStyle = 0

# Renumber the running index in column A to stay 0,1,2,...
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
